$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 25, shifting existing rows 25-27 down to 26-28
$ws.Rows.Item(25).Insert()

# Populate the new row with the monthly reports folder path entry
$ws.Range("A25").Value = "monthlyReportsFolderPath"
$ws.Range("B25").Value = "path\to\folder"
$ws.Range("B25").Style = "Normal"

# The last hyperlink (previously anchored at B27) needs to move to B28
# since the underlying row shifted down by one.
$ws.Range("B27").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B28"), "http://google.com/") | Out-Null

Write-Host ("Dimension: " + $ws.UsedRange.Address())
